$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells: AD1 = Wins, AE1 = Losses, AF1 = Ties
# Copy the style used by the existing header row (A1) so the new headers match.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill team record data (Wins=85, Losses=77, Ties=0) for every player row (2-49)
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 85   # AD
    $ws.Cells.Item($r, 31).Value = 77   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
